# Auto-generated Excel COM-interop edit script
# Applies numeric cell updates to the Kujata_Profits-style profit sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the commit diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 427.3
$ws.Range("I2").Value = 427.3
$ws.Range("K2").Value = 427.3
$ws.Range("M2").Value = -314.3
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("H19").Value = 423.92856
$ws.Range("I19").Value = 469.15384
$ws.Range("J19").Value = 384.73334
$ws.Range("K19").Value = 469.15384
$ws.Range("L19").Value = 384.73334
$ws.Range("M19").Value = -294.15384
$ws.Range("N19").Value = -734.73334
$ws.Range("H51").Value = 949.9
$ws.Range("I51").Value = 750
$ws.Range("J51").Value = 999.875
$ws.Range("K51").Value = 750
$ws.Range("L51").Value = 999.875
$ws.Range("M51").Value = -266
$ws.Range("N51").Value = -1967.875
$ws.Range("H62").Value = 15877873
$ws.Range("I62").Value = 18523852
$ws.Range("J62").Value = 2006
$ws.Range("K62").Value = 18523852
$ws.Range("L62").Value = 2006
$ws.Range("M62").Value = -18523228
$ws.Range("N62").Value = -3254
$ws.Range("H65").Value = 15877873
$ws.Range("I65").Value = 18523852
$ws.Range("J65").Value = 2006
$ws.Range("K65").Value = 92619260
$ws.Range("L65").Value = 10030
$ws.Range("M65").Value = -92616140
$ws.Range("N65").Value = -16270
$ws.Range("H74").Value = 2950.25
$ws.Range("I74").Value = 2001.5
$ws.Range("J74").Value = 3899
$ws.Range("K74").Value = 2001.5
$ws.Range("L74").Value = 3899
$ws.Range("M74").Value = -1065.5
$ws.Range("N74").Value = -5771
$ws.Range("H77").Value = 2950.25
$ws.Range("I77").Value = 2001.5
$ws.Range("J77").Value = 3899
$ws.Range("K77").Value = 10007.5
$ws.Range("L77").Value = 19495
$ws.Range("M77").Value = -5327.5
$ws.Range("N77").Value = -28855
$ws.Range("H132").Value = 8779982
$ws.Range("I132").Value = 9808929
$ws.Range("J132").Value = 33926.5
$ws.Range("K132").Value = 29426787
$ws.Range("L132").Value = 101779.5
$ws.Range("M132").Value = -29424257
$ws.Range("N132").Value = -106839.5
$ws.Range("H138").Value = 3050.0435
$ws.Range("J138").Value = 3122.1228
$ws.Range("L138").Value = 9366.368399999999
$ws.Range("N138").Value = -19646.3684

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18620.334
$ws.Range("I32").Value = 11072.046
$ws.Range("K32").Value = 11072.046
$ws.Range("M32").Value = -10785.046
$ws.Range("H74").Value = 1539.4147
$ws.Range("I74").Value = 857.8182
$ws.Range("J74").Value = 2328.6316
$ws.Range("K74").Value = 857.8182
$ws.Range("L74").Value = 2328.6316
$ws.Range("M74").Value = 16.18179999999995
$ws.Range("N74").Value = -4076.6316
$ws.Range("H77").Value = 1539.4147
$ws.Range("I77").Value = 857.8182
$ws.Range("J77").Value = 2328.6316
$ws.Range("K77").Value = 4289.091
$ws.Range("L77").Value = 11643.158
$ws.Range("M77").Value = 78.90899999999965
$ws.Range("N77").Value = -20379.158
$ws.Range("H122").Value = 4024.4
$ws.Range("I122").Value = 3413.875
$ws.Range("K122").Value = 10241.625
$ws.Range("M122").Value = -7791.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("M11").ClearContents()
$ws.Range("N11").ClearContents()
$ws.Range("H20").Value = 3194.2666
$ws.Range("I20").Value = 3086.8572
$ws.Range("J20").Value = 3288.25
$ws.Range("K20").Value = 3086.8572
$ws.Range("L20").Value = 3288.25
$ws.Range("M20").Value = -2839.8572
$ws.Range("N20").Value = -3782.25
$ws.Range("H132").Value = 39780
$ws.Range("J132").Value = 39780
$ws.Range("L132").Value = 39780
$ws.Range("N132").Value = -49900
$ws.Range("H134").Value = 2750
$ws.Range("I134").Value = 396.6
$ws.Range("J134").Value = 8633.5
$ws.Range("K134").Value = 1189.8
$ws.Range("L134").Value = 25900.5
$ws.Range("M134").Value = 1345.2
$ws.Range("N134").Value = -30970.5
$ws.Range("H140").Value = 64999.5
$ws.Range("J140").Value = 64999.5
$ws.Range("L140").Value = 64999.5
$ws.Range("N140").Value = -75359.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1735.0834
$ws.Range("I31").Value = 1708.1702
$ws.Range("J31").Value = 3000
$ws.Range("K31").Value = 1708.1702
$ws.Range("L31").Value = 3000
$ws.Range("M31").Value = -1413.1702
$ws.Range("N31").Value = -3590
$ws.Range("H34").Value = 1735.0834
$ws.Range("I34").Value = 1708.1702
$ws.Range("J34").Value = 3000
$ws.Range("K34").Value = 1708.1702
$ws.Range("L34").Value = 3000
$ws.Range("M34").Value = -1506.1702
$ws.Range("N34").Value = -3404
$ws.Range("H62").Value = 20002420
$ws.Range("I62").Value = 2600
$ws.Range("J62").Value = 66668668
$ws.Range("K62").Value = 2600
$ws.Range("L62").Value = 66668668
$ws.Range("M62").Value = -1976
$ws.Range("N62").Value = -66669916
$ws.Range("H65").Value = 20002420
$ws.Range("I65").Value = 2600
$ws.Range("J65").Value = 66668668
$ws.Range("K65").Value = 13000
$ws.Range("L65").Value = 333343340
$ws.Range("M65").Value = -9880
$ws.Range("N65").Value = -333349580
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()
$ws.Range("H132").Value = 1719.8837
$ws.Range("I132").Value = 1403.0646
$ws.Range("J132").Value = 2538.3333
$ws.Range("K132").Value = 4209.1938
$ws.Range("L132").Value = 7614.999899999999
$ws.Range("M132").Value = -1679.1938
$ws.Range("N132").Value = -12674.9999
$ws.Range("H141").Value = 424570.4
$ws.Range("J141").Value = 424570.4
$ws.Range("L141").Value = 424570.4
$ws.Range("N141").Value = -434930.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H49").Value = 5668
$ws.Range("I49").Value = 5000
$ws.Range("J49").Value = 6002
$ws.Range("K49").Value = 15000
$ws.Range("L49").Value = 18006
$ws.Range("M49").Value = -14844
$ws.Range("N49").Value = -18318
$ws.Range("H56").Value = 5553
$ws.Range("I56").Value = 5553
$ws.Range("K56").Value = 5553
$ws.Range("M56").Value = -5023
$ws.Range("H107").Value = 7257.25
$ws.Range("J107").Value = 18746.334
$ws.Range("L107").Value = 56239.00199999999
$ws.Range("N107").Value = -60079.00199999999
$ws.Range("H123").Value = 2316.25
$ws.Range("I123").Value = 1906
$ws.Range("K123").Value = 5718
$ws.Range("M123").Value = -3268
$ws.Range("H131").Value = 22256496
$ws.Range("J131").Value = 43957.027
$ws.Range("L131").Value = 131871.081
$ws.Range("N131").Value = -141951.081
$ws.Range("H136").Value = 2664.5
$ws.Range("I136").Value = 1025
$ws.Range("J136").Value = 4959.8
$ws.Range("K136").Value = 3075
$ws.Range("L136").Value = 14879.4
$ws.Range("M136").Value = 2025
$ws.Range("N136").Value = -25079.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 22503662
$ws.Range("I70").Value = 20837068
$ws.Range("J70").Value = 25003552
$ws.Range("K70").Value = 20837068
$ws.Range("L70").Value = 25003552
$ws.Range("M70").Value = -20836798
$ws.Range("N70").Value = -25004092
$ws.Range("H73").Value = 22503662
$ws.Range("I73").Value = 20837068
$ws.Range("J73").Value = 25003552
$ws.Range("K73").Value = 20837068
$ws.Range("L73").Value = 25003552
$ws.Range("M73").Value = -20836132
$ws.Range("N73").Value = -25005424
$ws.Range("H122").Value = 1450.1666
$ws.Range("I122").Value = 1131.7693
$ws.Range("J122").Value = 2278
$ws.Range("K122").Value = 3395.3079
$ws.Range("L122").Value = 6834
$ws.Range("M122").Value = -945.3078999999998
$ws.Range("N122").Value = -11734
$ws.Range("H132").Value = 8109.364
$ws.Range("I132").Value = 9170.235000000001
$ws.Range("J132").Value = 4502.4
$ws.Range("K132").Value = 27510.705
$ws.Range("L132").Value = 13507.2
$ws.Range("M132").Value = -24980.705
$ws.Range("N132").Value = -18567.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 331.06668
$ws.Range("I55").Value = 199.5
$ws.Range("J55").Value = 857.3333
$ws.Range("K55").Value = 199.5
$ws.Range("L55").Value = 857.3333
$ws.Range("M55").Value = -26.5
$ws.Range("N55").Value = -1203.3333
$ws.Range("H122").Value = 47224456
$ws.Range("J122").Value = 3351
$ws.Range("L122").Value = 10053
$ws.Range("N122").Value = -14953
$ws.Range("H132").Value = 3554.0908
$ws.Range("I132").Value = 3019.8
$ws.Range("J132").Value = 3999.3333
$ws.Range("K132").Value = 9059.400000000001
$ws.Range("L132").Value = 11997.9999
$ws.Range("M132").Value = -6529.400000000001
$ws.Range("N132").Value = -17057.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1781.375
$ws.Range("J81").Value = 1813.8966
$ws.Range("L81").Value = 3627.7932
$ws.Range("N81").Value = -5749.7932
$ws.Range("H84").Value = 1781.375
$ws.Range("J84").Value = 1813.8966
$ws.Range("L84").Value = 18138.966
$ws.Range("N84").Value = -28746.966
$ws.Range("H122").Value = 15627320
$ws.Range("I122").Value = 17859494
$ws.Range("J122").Value = 2100
$ws.Range("K122").Value = 53578482
$ws.Range("L122").Value = 6300
$ws.Range("M122").Value = -53576032
$ws.Range("N122").Value = -11200
$ws.Range("H132").Value = 4507.6978
$ws.Range("I132").Value = 4299.3
$ws.Range("J132").Value = 4988.615
$ws.Range("K132").Value = 12897.9
$ws.Range("L132").Value = 14965.845
$ws.Range("M132").Value = -10367.9
$ws.Range("N132").Value = -20025.845
$ws.Range("H136").Value = 2100.2942
$ws.Range("I136").Value = 1907.5
$ws.Range("K136").Value = 5722.5
$ws.Range("M136").Value = -3172.5
$ws.Range("H140").Value = 34002.25
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 34002.25
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 34002.25
$ws.Range("M140").ClearContents()
$ws.Range("N140").Value = -44362.25

Write-Output "Updated cells: sets=258 clears=5"
